# Add a "2022-Q1" holdings sheet (mirrors the 2021-Q4 layout) right before the
# "总计" summary sheet, and prepend a matching "2022-Q1" row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned immediately before "总计".
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"
# Match the <outlinePr summaryBelow="1" summaryRight="1"/> that every other
# quarterly sheet carries in <sheetPr>.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Clone the full header/data block (values + styles) from the 2021-Q4 sheet so
# fonts/borders/alignment match the other quarterly sheets exactly.
$templateSheet.Range("A1:H29").Copy($newSheet.Range("A1"))
# The new sheet only needs 25 data rows (26 total incl. header); drop the tail.
$newSheet.Range("A27:H29").Delete()
# Header row has no A1 cell in the template sheets (B1:H1 only) - the Copy()
# above materialised an empty A1 shell, so drop it to match.
$newSheet.Range("A1").ClearContents()

function Set-TextCell($cell, [string]$text) {
    # Force the cell to remain text (Excel auto-coerces numeric-looking
    # strings to numbers on assignment) using the same leading-apostrophe
    # trick a user would type, then drop the resulting quote-prefix style so
    # no stray formatting is left behind.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$fundRows = @(
    @("501207", "华夏创新未来18个月封闭运作混合A", "67.75", "84.92", "3.14", "2.1274", 9),
    @("000031", "华夏复兴混合", "27.37", "89.15", "7.49", "2.0500", 3),
    @("910007", "东方红启元三年持有期混合A", "69.33", "74.32", "1.79", "1.2410", 10),
    @("007887", "东方红启元三年持有期混合B", "59.71", "74.32", "1.79", "1.0688", 10),
    @("340008", "兴全有机增长混合", "30.22", "75.70", "3.44", "1.0396", 8),
    @("009951", "广发稳健回报混合A", "70.64", "43.16", "1.20", "0.8477", 10),
    @("007349", "华夏科技创新混合A", "14.63", "89.16", "4.26", "0.6232", 9),
    @("011138", "广发聚鸿六个月持有期混合A", "12.07", "93.49", "5.04", "0.6083", 6),
    @("011140", "广发聚鸿六个月持有期混合E", "12.07", "93.49", "5.04", "0.6083", 6),
    @("010106", "华夏核心科技6个月定期开放混合A", "8.53", "79.73", "6.10", "0.5203", 4),
    @("501070", "广发睿阳三年定期开放混合", "7.06", "50.14", "6.58", "0.4645", 1),
    @("162720", "广发创业板两年定期开放混合", "8.96", "93.33", "5.02", "0.4498", 5),
    @("010518", "华夏先锋科技一年定期开放混合A", "8.72", "88.08", "5.05", "0.4404", 7),
    @("169107", "东方红恒阳五年定期开放混合", "22.26", "78.45", "1.93", "0.4296", 8),
    @("002124", "广发新兴产业精选灵活配置混合", "11.11", "91.11", "3.71", "0.4122", 6),
    @("002472", "光大保德信先进服务业灵活配置混合", "6.67", "89.61", "3.98", "0.2655", 9),
    @("519673", "银河康乐股票", "2.31", "92.35", "6.04", "0.1395", 3),
    @("009952", "广发稳健回报混合C", "11.53", "43.16", "1.20", "0.1384", 10),
    @("008234", "光大保德信消费主题股票", "2.26", "91.66", "5.35", "0.1209", 6),
    @("010519", "华夏先锋科技一年定期开放混合C", "2.12", "88.08", "5.05", "0.1071", 7),
    @("010107", "华夏核心科技6个月定期开放混合C", "1.32", "79.73", "6.10", "0.0805", 4),
    @("007350", "华夏科技创新混合C", "1.08", "89.16", "4.26", "0.0460", 9),
    @("011139", "广发聚鸿六个月持有期混合C", "0.64", "93.49", "5.04", "0.0323", 6),
    @("005027", "光大保德信多策略优选一年定期开放灵活配置混合", "0.52", "37.88", "1.99", "0.0103", 10),
    @("005166", "嘉实润和量化6个月定期开放混合", "0.25", "28.26", "0.66", "0.0016", 3)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i
    Set-TextCell $newSheet.Cells.Item($r, 2) $row[0]
    Set-TextCell $newSheet.Cells.Item($r, 3) $row[1]
    Set-TextCell $newSheet.Cells.Item($r, 4) $row[2]
    Set-TextCell $newSheet.Cells.Item($r, 5) $row[3]
    Set-TextCell $newSheet.Cells.Item($r, 6) $row[4]
    Set-TextCell $newSheet.Cells.Item($r, 7) $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the rest
#    down one row and renumbering the A-column sequence index.
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")

$existingSummaryRows = @()
for ($r = 2; $r -le 6; $r++) {
    $existingSummaryRows += ,@(
        $summarySheet.Cells.Item($r, 2).Value(),
        $summarySheet.Cells.Item($r, 3).Value(),
        $summarySheet.Cells.Item($r, 4).Value()
    )
}

$styleSource = $summarySheet.Range("A2")

$summaryRows = New-Object System.Collections.ArrayList
[void]$summaryRows.Add(@("2022-Q1", 25, 13.87))
foreach ($row in $existingSummaryRows) { [void]$summaryRows.Add($row) }

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]

    $styleSource.Copy($summarySheet.Cells.Item($r, 1))
    $summarySheet.Cells.Item($r, 1).Value = $i
    $summarySheet.Cells.Item($r, 2).Value = $row[0]
    $summarySheet.Cells.Item($r, 3).Value = $row[1]
    $summarySheet.Cells.Item($r, 4).Value = $row[2]
}
